# Guest user 3d for drybar UK
# Inserts a new row (row 13) into the "DataSet" sheet carrying a 3-D Secure
# guest-checkout card row, shifting all subsequent rows down by one.
# Because this engine's Hyperlinks collection does not automatically
# re-target its cached Range when rows are inserted, we capture the
# existing hyperlinks first, remove them, perform the insert, and then
# recreate them pointing at their (possibly shifted) destinations.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSet")

# ---------------------------------------------------------------------
# 1. Record the current hyperlinks (address/display/row/col) so they can
#    be rebuilt after the row insert shifts everything below row 13.
# ---------------------------------------------------------------------
$hlCount = $ws.Hyperlinks.Count
$hlRow = @()
$hlCol = @()
$hlAddress = @()
$hlDisplay = @()

for ($i = 1; $i -le $hlCount; $i++) {
    $h = $ws.Hyperlinks.Item($i)
    $hlRow += $h.Range.Row
    $hlCol += $h.Range.Column
}

# Known external targets for each existing hyperlink, in sheet order
# (row, col, target-address, display-text-or-null).
$hlData = @(
    @(17, 11, "mailto:avayugundla@helenoftroy.com", $null),
    @(18, 11, "mailto:avayugundla@helenoftroy.com", $null),
    @(2, 4, "mailto:Lotuswave@123", $null),
    @(2, 5, "mailto:Lotuswave@123", $null),
    @(2, 3, "mailto:avayugundla@helenoftroy.com", $null),
    @(2, 2, "mailto:avayugundla@helenoftroy.com", $null),
    @(21, 11, "mailto:avayugundla@helenoftroy.com", $null),
    @(23, 11, "mailto:Paypal-buyer@hydroflask.com", $null),
    @(25, 11, "mailto:avayugundla@helenoftroy.com", $null),
    @(34, 2, "mailto:hydroflaskemea978@gmail.com", $null),
    @(34, 4, "mailto:Lotus@123", $null),
    @(34, 5, "mailto:Lotus@123", $null),
    @(34, 11, "mailto:hydroflaskemea978@gmail.com", $null),
    @(39, 11, "mailto:qatesting.lotuswave@gmail.com", $null),
    @(40, 2, "mailto:mmarella@helenoftroy.com", $null),
    @(40, 4, "mailto:Lotus@123", $null),
    @(40, 5, "mailto:Lotus@123", $null),
    @(40, 3, "mailto:mmarella@helenoftroy.com", $null),
    @(8, 6, "mailto:!#@", $null),
    @(8, 4, "mailto:Lotus@1", $null),
    @(8, 5, "mailto:Lotus@1235", $null),
    @(42, 2, "mailto:hydroflaskemea978+7@gmail.com", "mailto:hydroflaskemea978+7@gmail.com"),
    @(42, 3, "mailto:hydroflaskemea978+7@gmail.com", "mailto:hydroflaskemea978+7@gmail.com"),
    @(42, 4, "mailto:Lotus@123", "mailto:Lotus@123"),
    @(42, 5, "mailto:Lotus@123", "mailto:Lotus@123"),
    @(42, 11, "mailto:hydroflaskemea978+7@gmail.com", "mailto:hydroflaskemea978+7@gmail.com"),
    @(43, 2, "mailto:hydroflaskemea978+8@gmail.com", "mailto:hydroflaskemea978+8@gmail.com"),
    @(43, 4, "mailto:Lotus@123", "mailto:Lotus@123"),
    @(43, 5, "mailto:Lotuswave@1234", "mailto:Lotuswave@1234"),
    @(43, 11, "mailto:avayugundla@helenoftroy.com", "mailto:avayugundla@helenoftroy.com"),
    @(43, 12, "mailto:hydroflaskemea978+2@gmail.com", "mailto:hydroflaskemea978+2@gmail.com"),
    @(46, 11, "mailto:avayugundla@helenoftroy.com", $null)
)

# Remove all existing hyperlinks up front (individual Delete() calls are
# unreliable here, but clearing the whole collection works).
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------
# 2. Insert the new row at position 13; everything at/after row 13 moves
#    down to make room (rows 13-46 become rows 14-47).
# ---------------------------------------------------------------------
$ws.Range("A13").EntireRow.Insert()

# ---------------------------------------------------------------------
# 3. Populate the newly inserted guest/3-D Secure card row.
# ---------------------------------------------------------------------
$ws.Range("A13").Value2 = "3d_Secure"
$ws.Range("U13").Value2 = "'4000000000003220"
$ws.Range("V13").Value2 = "'06/29"
$ws.Range("W13").Value2 = 123
$ws.Range("AJ13").Value2 = [char]0x00A3

# ---------------------------------------------------------------------
# 4. Recreate the hyperlinks, shifting any that sat at row >= 13 down by
#    one row to follow the cells they were attached to.
# ---------------------------------------------------------------------
foreach ($entry in $hlData) {
    $origRow = $entry[0]
    $col = $entry[1]
    $address = $entry[2]
    $display = $entry[3]

    $newRow = $origRow
    if ($origRow -ge 13) {
        $newRow = $origRow + 1
    }

    $target = $ws.Cells.Item($newRow, $col)
    if ($display) {
        $ws.Hyperlinks.Add($target, $address, [System.Type]::Missing, [System.Type]::Missing, $display)
    } else {
        $ws.Hyperlinks.Add($target, $address)
    }
}

# ---------------------------------------------------------------------
# 5. Update the sheet view: scroll/selection moved from O28 to AJ17, and
#    the visible window now starts around column AG.
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 33
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("AJ17").Select()
